{"js": "// The worksheet's single table has 20 rows x 5 columns; only rows\n// 0, 4, 8, 12, 16 hold the \"NN\u00f7N=\" practice prompts (the rows in\n// between are blank answer rows). Update each of those 25 cells, in\n// document order, to the new values from the commit. We target cells\n// by position (row/column) rather than by text search, because some\n// of the new values collide with other cells' old values (e.g. a\n// cell's old text \"23\u00f79=\" becomes \"14\u00f77=\", while another cell's old\n// text already was \"14\u00f77=\" and must become \"34\u00f75=\") \u2014 a plain\n// sequential search-and-replace would incorrectly re-match already\n// updated text.\nconst newValuesByRow = [\n  [\"17\u00f78=\", \"30\u00f75=\", \"14\u00f77=\", \"31\u00f74=\", \"89\u00f77=\"],\n  [\"59\u00f79=\", \"48\u00f75=\", \"22\u00f79=\", \"34\u00f75=\", \"65\u00f76=\"],\n  [\"52\u00f78=\", \"58\u00f76=\", \"32\u00f76=\", \"27\u00f72=\", \"86\u00f72=\"],\n  [\"57\u00f75=\", \"67\u00f72=\", \"65\u00f79=\", \"41\u00f78=\", \"42\u00f75=\"],\n  [\"33\u00f78=\", \"40\u00f78=\", \"77\u00f78=\", \"92\u00f73=\", \"79\u00f77=\"],\n];\nconst contentRowIndexes = [0, 4, 8, 12, 16];\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let ri = 0; ri < contentRowIndexes.length; ri++) {\n  const row = rows.items[contentRowIndexes[ri]];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let ci = 0; ci < cells.items.length; ci++) {\n    const cellBody = cells.items[ci].body;\n    const paragraphs = cellBody.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    const paragraphRange = paragraphs.items[0].getRange();\n    paragraphRange.insertText(newValuesByRow[ri][ci], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The worksheet's single table has 20 rows x 5 columns; only rows\n# 1, 5, 9, 13, 17 (1-based) hold the \"NN\u00f7N=\" practice prompts (the\n# rows in between are blank answer rows). Update each of those 25\n# cells, in document order, to the new values from the commit.\n#\n# We target cells by (row, column) position via $table.Cell(r, c)\n# rather than Find/Replace on text, because some of the new values\n# collide with other cells' old values (e.g. one cell's old text\n# \"23\u00f79=\" becomes \"14\u00f77=\", while a different cell's old text already\n# was \"14\u00f77=\" and must become \"34\u00f75=\"). A global Find/Replace would\n# incorrectly re-match text that a previous replacement just wrote.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$contentRows = @(1, 5, 9, 13, 17)\n$newValues = @(\n    @(\"17\u00f78=\", \"30\u00f75=\", \"14\u00f77=\", \"31\u00f74=\", \"89\u00f77=\"),\n    @(\"59\u00f79=\", \"48\u00f75=\", \"22\u00f79=\", \"34\u00f75=\", \"65\u00f76=\"),\n    @(\"52\u00f78=\", \"58\u00f76=\", \"32\u00f76=\", \"27\u00f72=\", \"86\u00f72=\"),\n    @(\"57\u00f75=\", \"67\u00f72=\", \"65\u00f79=\", \"41\u00f78=\", \"42\u00f75=\"),\n    @(\"33\u00f78=\", \"40\u00f78=\", \"77\u00f78=\", \"92\u00f73=\", \"79\u00f77=\")\n)\n\nfor ($ri = 0; $ri -lt $contentRows.Length; $ri++) {\n    $row = $contentRows[$ri]\n    for ($ci = 0; $ci -lt 5; $ci++) {\n        $col = $ci + 1\n        $cell = $table.Cell($row, $col)\n        $cell.Range.Text = $newValues[$ri][$ci]\n    }\n}\n"}
